# Add a new 13th slide ("My aliases") using the "Title and Content" layout
# (same layout as the existing "Some references" slide, slideLayout2.xml /
# ppLayout 2), appended at the end of the slide list.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = 'My aliases'

$line1 = 'alias.co=checkout'
$line2 = 'alias.l=log --graph --pretty=format:''%Cred%h%Creset -%C(auto)%d%Creset %s %Cgreen(%cr) %C(bold blue)<%an>%Creset'' --abbrev-commit --date=relative'
$line3 = 'alias.cod=checkout paulmey/dev'
$line4 = 'alias.s=status -s -b'
$line5 = 'alias.fa=fetch --all'

$body = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5

$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = $body
